# Rename transcript speaker tags: replace "T/R2" with "T" in column D
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)   # Column D
    if ($cell.Value2 -eq "T/R2") {
        $cell.Value = "T"
    }
}
